$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111272062
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("Q2").Value = 473156.3705774212
$ws.Range("R2").Value = 6863531.269191674
$ws.Range("A3").Value = 111271821
$ws.Range("Q3").Value = 473140.3516782348
$ws.Range("R3").Value = 6863595.022241795
$ws.Range("AC3").Value = "Tre blommande."
$ws.Range("A4").Value = 111271296
$ws.Range("B4").Value = 78578
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 473220.1559155915
$ws.Range("R4").Value = 6863539.25170773
$ws.Range("A5").Value = 111270939
$ws.Range("Q5").Value = 473229.5908188519
$ws.Range("R5").Value = 6863658.889402787
$ws.Range("A6").Value = 111271055
$ws.Range("B6").Value = 78579
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 2081
$ws.Range("F6").Value = "Skrovellav"
$ws.Range("G6").Value = "Lobaria scrobiculata"
$ws.Range("H6").Value = "(Scop.) DC."
$ws.Range("Q6").Value = 473238.8676645419
$ws.Range("R6").Value = 6863638.079474191
$ws.Range("AC6").Value = $null
$ws.Range("A7").Value = 111271923
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("Q7").Value = 473118.5439814709
$ws.Range("R7").Value = 6863582.939962601
$ws.Range("AC7").Value = "Tre blommande."
$ws.Range("A8").Value = 111271141
$ws.Range("B8").Value = 78578
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = "Lunglav"
$ws.Range("G8").Value = "Lobaria pulmonaria"
$ws.Range("H8").Value = "(L.) Hoffm."
$ws.Range("A9").Value = 111270747
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("Q9").Value = 473194.7999623233
$ws.Range("R9").Value = 6863736.454484907
$ws.Range("AC9").Value = "Sex blommande."
$ws.Range("A10").Value = 111271061
$ws.Range("Q10").Value = 473238.8676645419
$ws.Range("R10").Value = 6863638.079474191
$ws.Range("AC10").Value = $null
$ws.Range("A11").Value = 111270784
$ws.Range("B11").Value = 78578
$ws.Range("E11").Value = 6458
$ws.Range("F11").Value = "Lunglav"
$ws.Range("G11").Value = "Lobaria pulmonaria"
$ws.Range("H11").Value = "(L.) Hoffm."
$ws.Range("Q11").Value = 473239.9383552746
$ws.Range("R11").Value = 6863714.420922431
$ws.Range("A12").Value = 111271382
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 473167.6377000402
$ws.Range("R12").Value = 6863583.496200636
$ws.Range("AC12").Value = "Tre blommande."
$ws.Range("A13").Value = 111271588
$ws.Range("B13").Value = 78578
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6458
$ws.Range("F13").Value = "Lunglav"
$ws.Range("G13").Value = "Lobaria pulmonaria"
$ws.Range("H13").Value = "(L.) Hoffm."
$ws.Range("Q13").Value = 473140.3516782348
$ws.Range("R13").Value = 6863595.022241795
$ws.Range("A14").Value = 111272343
$ws.Range("B14").Value = 78578
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6458
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."
$ws.Range("Q14").Value = 473387.8703240218
$ws.Range("R14").Value = 6863558.206130736
$ws.Range("A15").Value = 111272375
$ws.Range("B15").Value = 96251
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 220093
$ws.Range("F15").Value = "Korallrot"
$ws.Range("G15").Value = "Corallorhiza trifida"
$ws.Range("H15").Value = "Châtel."
$ws.Range("Q15").Value = 473400.7315261344
$ws.Range("R15").Value = 6863573.187783281
$ws.Range("A16").Value = 111270596
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("Q16").Value = 473184.8241620373
$ws.Range("R16").Value = 6863788.37406126
$ws.Range("AC16").Value = "Fem blommande."
$ws.Range("A17").Value = 111271309
$ws.Range("B17").Value = 78579
$ws.Range("E17").Value = 2081
$ws.Range("F17").Value = "Skrovellav"
$ws.Range("G17").Value = "Lobaria scrobiculata"
$ws.Range("H17").Value = "(Scop.) DC."
$ws.Range("Q17").Value = 473221.4734201821
$ws.Range("R17").Value = 6863586.84377678
$ws.Range("A18").Value = 111271176
$ws.Range("B18").Value = 78579
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 2081
$ws.Range("F18").Value = "Skrovellav"
$ws.Range("G18").Value = "Lobaria scrobiculata"
$ws.Range("H18").Value = "(Scop.) DC."
$ws.Range("Q18").Value = 473227.9160841404
$ws.Range("R18").Value = 6863625.911539786
$ws.Range("AC18").Value = $null
$ws.Range("A19").Value = 111271029
$ws.Range("B19").Value = 78579
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 2081
$ws.Range("F19").Value = "Skrovellav"
$ws.Range("G19").Value = "Lobaria scrobiculata"
$ws.Range("H19").Value = "(Scop.) DC."
$ws.Range("Q19").Value = 473229.5908188519
$ws.Range("R19").Value = 6863658.889402787
$ws.Range("AC19").Value = $null
$ws.Range("A20").Value = 111272292
$ws.Range("B20").Value = 78578
$ws.Range("E20").Value = 6458
$ws.Range("F20").Value = "Lunglav"
$ws.Range("G20").Value = "Lobaria pulmonaria"
$ws.Range("H20").Value = "(L.) Hoffm."
$ws.Range("Q20").Value = 473321.1690919191
$ws.Range("R20").Value = 6863539.403128584
$ws.Range("A21").Value = 111270755
$ws.Range("B21").Value = 78578
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("Q21").Value = 473194.7999623233
$ws.Range("R21").Value = 6863736.454484907
$ws.Range("AC21").Value = $null
$ws.Range("A22").Value = 111270559
$ws.Range("B22").Value = 96348
$ws.Range("Q22").Value = 473167.8634183492
$ws.Range("R22").Value = 6863792.277629613
$ws.Range("AC22").Value = "Två blommande."
